$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate column I (Image Name values) with a CONCATENATE formula that
# builds "<ItemName>.png" from column B, for every data row (2-57).
for ($r = 2; $r -le 57; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.Formula = "=CONCATENATE(B" + $r + ",`".png`")"
}

# Resize column I (bestFit-style width) to fit the new content.
$ws.Columns.Item(9).ColumnWidth = 14.5

# Move / update the active selection as recorded in the saved view state.
$ws.Range("K12").Select()
